$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.820.59'
$ws.Cells.Item(2, 5).Value = '  +0.84%  '

$ws.Cells.Item(3, 4).Value = '2.492.41'
$ws.Cells.Item(3, 5).Value = '  -0.05%  '

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '586.86'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.23%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '177.08'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.73%  '

$ws.Cells.Item(7, 5).Value = '  -0.01%  '

$ws.Cells.Item(8, 5).Value = '  +0.16%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.141'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +3.57%  '

$ws.Cells.Item(10, 5).Value = '  -0.21%  '

$ws.Cells.Item(11, 5).Value = '  +1.76%  '

$ws.Cells.Item(12, 5).Value = '  -0.02%  '

$ws.Cells.Item(13, 4).Value = '2.948.64'
$ws.Cells.Item(13, 5).Value = '  +0.99%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.66'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.54%  '

$ws.Cells.Item(15, 4).Value = '67.715.42'
$ws.Cells.Item(15, 5).Value = '  +0.80%  '

$ws.Cells.Item(16, 5).Value = '  +0.55%  '

$ws.Cells.Item(17, 4).Value = '2.492.32'
$ws.Cells.Item(17, 5).Value = '  -0.19%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.97'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.89%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.48'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.32%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '351.22'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.07%  '

$ws.Cells.Item(21, 5).Value = '  +2.01%  '

$ws.Cells.Item(22, 5).Value = '  +0.04%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '70.83'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +3.14%  '

$ws.Cells.Item(24, 5).Value = '  +0.51%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.74'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -3.00%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.11'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.73%  '

$ws.Cells.Item(27, 4).Value = '2.620.32'
$ws.Cells.Item(27, 5).Value = '  +0.01%  '

$ws.Cells.Item(28, 5).Value = '  +0.07%  '

$ws.Cells.Item(29, 4).Value = '0.0₃0905'
$ws.Cells.Item(29, 5).Value = '  +0.08%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '505.62'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.31%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.82'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.16%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.26'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.64%  '

$ws.Cells.Item(33, 5).Value = '  -0.09%  '

$ws.Cells.Item(34, 5).Value = '  +0.01%  '

$ws.Cells.Item(35, 5).Value = '  +2.92%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '163.31'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.90%  '

$ws.Cells.Item(37, 5).Value = '  -0.25%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.33'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.33'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.36%  '

$ws.Cells.Item(40, 5).Value = '  -0.05%  '

$ws.Cells.Item(41, 5).Value = '  +3.01%  '

$ws.Cells.Item(42, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.329'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.14%  '

$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.86'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.23%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.42'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +2.30%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '144.62'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.85%  '

$ws.Cells.Item(46, 5).Value = '  +1.75%  '

$ws.Cells.Item(47, 5).Value = '  -0.19%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.58'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.60%  '

$ws.Cells.Item(50, 5).Value = '  +0.11%  '

$ws.Cells.Item(51, 5).Value = '  +0.39%  '
